$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# "Generate Report for Handoff"
#
# The 65a68d0c... file has moved from "Handed back: in sync with en-US"
# to "Ready for handoff" (new handoff timestamps), and the d5deb3a4...
# file's row (which was mid-handback) is removed entirely from every
# sheet since it's no longer tracked in this report. The row that used
# to be below it (.localization-config) shifts up to take its place.
# -----------------------------------------------------------------------

# ============================== Overview ================================
$ws1 = $wb.Worksheets.Item("Overview")

# Remove the d5deb3a4-...md row; .localization-config row shifts from 4 -> 3
$ws1.Rows.Item(3).Delete()

# Status text for the remaining tracked file
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# Hyperlinks need to be rebuilt: this engine does not re-anchor the
# <hyperlinks> table when rows shift, so clear and re-add them in the
# final layout.
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0d41aa794f7c40cb62fa00651e6f26cde25f00b2/e2e/65a68d0c-23e2-493e-9abc-d1df473efd8c.md", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0d41aa794f7c40cb62fa00651e6f26cde25f00b2/.localization-config", "", "", ".localization-config")

# ================================ zh-cn ==================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(3).Delete()

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-03-10 11:58:00"

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0d41aa794f7c40cb62fa00651e6f26cde25f00b2/e2e/65a68d0c-23e2-493e-9abc-d1df473efd8c.md", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e2244a6e832bb2f19a8921acecef89242a9234d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.zh-cn.xlf", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a035ee8e5e1cadc2a919770f29c60dfa458a12b8/e2e/65a68d0c-23e2-493e-9abc-d1df473efd8c.md", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/149c46bd3c21bd65d6afd077ee5fe889500300c6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.zh-cn.xlf", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0d41aa794f7c40cb62fa00651e6f26cde25f00b2/.localization-config", "", "", ".localization-config")

# ================================ de-de ==================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(3).Delete()

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-03-10 11:58:04"

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0d41aa794f7c40cb62fa00651e6f26cde25f00b2/e2e/65a68d0c-23e2-493e-9abc-d1df473efd8c.md", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2ef4754d842fe0d4a0252b693f95ca6486e3dc13/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.de-de.xlf", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c5ae4dfe7d7fdd222ffc8ca1471e8d04cecb1bf6/e2e/65a68d0c-23e2-493e-9abc-d1df473efd8c.md", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d63149cd9c1c5aa8c6e7dc5623e8afd6c17b37ad/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.de-de.xlf", "", "", "65a68d0c-23e2-493e-9abc-d1df473efd8c.8d603966d2b8bc9a8d210f5aee3ac6ff558d6f1e.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0d41aa794f7c40cb62fa00651e6f26cde25f00b2/.localization-config", "", "", ".localization-config")

Write-Output "Report regenerated for handoff."
